$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$wsMeta = $wb.Worksheets.Item("Metadata")

# Publisher (row 9) value translated from German to English
$wsMeta.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"

# Contact (row 10) value translated from German to English
$wsMeta.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Description (row 12) - new value added
$wsMeta.Range("B12").Value = "consent states - subset CONSENT documents"
